# Apply the "last report 25-12-24" updates to the Product Requisition workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Data updates -----------------------------------------------------
# Row 9: set quantity (C9) so the Total Amount (E9) recalculates to 7050
$ws.Range("C9").Value = 15

# Row 14: clear quantity (C14) so Total Amount (E14) recalculates to 0
$ws.Range("C14").ClearContents()

# Row 31: reduce quantity (C31) from 10000 to 2000
$ws.Range("C31").Value = 2000

# Row 32: clear quantity (C32) so Total Amount (E32) recalculates to 0
$ws.Range("C32").ClearContents()

# Row 43: reduce quantity (C43) from 243909 to 112551
$ws.Range("C43").Value = 112551

# Let Excel recalculate dependent formulas (E column + the E46 grand total).
$excel.Calculate()

# --- View state -----------------------------------------------------
# Move the visible window up a bit and land the selection on C44.
$ws.Range("A22").Select() | Out-Null
$ws.Range("C44").Select()
